# Update population margins: restrict to SBV residents 18+ (drop under-18 buckets,
# collapse 80-84 / 85+ into a single "80 +" bucket), and recompute the totals/proportions
# to match the new (COM-derived) 2020-vintage population counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the now-unused bottom rows (old rows 30 and 32:37 -- everything
#    below the "male / 80 +" row). Delete the trailing block first (A32:H37)
#    so row numbers above stay put, then wipe row 30 completely (so it drops
#    out of the sheet), leaving row 31 as a bare, number-formatted spacer
#    cell in column D (mirrors the source row that still carries the old
#    "0.0000" style but no value).
# ---------------------------------------------------------------------------
$ws.Range("A32:H37").Delete() | Out-Null
$ws.Range("A30:H30").Clear() | Out-Null
$ws.Range("A31:C31").ClearContents() | Out-Null
$ws.Range("D31").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 2. Age-group labels (column A) and sex labels (column B).
#    Female block: rows 2-14 = 15-19 .. 75-79, row 15 = "80 +"
#    Male block:   rows 16-28 = 15-19 .. 75-79, row 29 = "80 +"
# ---------------------------------------------------------------------------
$ageLabels = @("15 to 19 years","20 to 24 years","25 to 29 years","30 to 34 years","35 to 39 years","40 to 44 years","45 to 49 years","50 to 54 years","55 to 59 years","60 to 64 years ","65 to 69 years","70 to 74 years","75 to 79 years")

for ($i = 0; $i -lt $ageLabels.Length; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $ageLabels[$i]
    $ws.Cells.Item(2 + $i, 2).Value = "female"
}
$ws.Cells.Item(15, 1).Value = "80 +"
$ws.Cells.Item(15, 2).Value = "female"

for ($i = 0; $i -lt $ageLabels.Length; $i++) {
    $ws.Cells.Item(16 + $i, 1).Value = $ageLabels[$i]
    $ws.Cells.Item(16 + $i, 2).Value = "male"
}
$ws.Cells.Item(29, 1).Value = "80 +"
$ws.Cells.Item(29, 2).Value = "male"

# ---------------------------------------------------------------------------
# 3. Raw population counts (column C).
# ---------------------------------------------------------------------------
$femaleCounts = @{
    3  = 49829
    4  = 44444
    5  = 38862
    6  = 37618
    7  = 39026
    8  = 39220
    9  = 39659
    10 = 42425
    11 = 40867
    12 = 33498
    13 = 28663
    14 = 18795
}
$ws.Cells.Item(2, 3).Formula = "=48510*(1/5)"
foreach ($row in $femaleCounts.Keys) {
    $ws.Cells.Item($row, 3).Value = $femaleCounts[$row]
}
$ws.Cells.Item(15, 3).Formula = "=13702+19247"

$maleCounts = @{
    17 = 52703
    18 = 46176
    19 = 43673
    20 = 39211
    21 = 39665
    22 = 38286
    23 = 40016
    24 = 42179
    25 = 37308
    26 = 28370
    27 = 26908
    28 = 16073
}
$ws.Cells.Item(16, 3).Formula = "=47562*(1/5)"
foreach ($row in $maleCounts.Keys) {
    $ws.Cells.Item($row, 3).Value = $maleCounts[$row]
}
$ws.Cells.Item(29, 3).Formula = "=11073+10624"

# ---------------------------------------------------------------------------
# 4. Total population denominator (H2) and proportion formulas (column D).
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 8).Formula = "=481777.4+495557"

for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 4).Formula = "=C$row/`$H`$2"
}

# ---------------------------------------------------------------------------
# 5. Misc cosmetic bits that shifted along with the data edit.
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 71
$ws.Range("D31").Select() | Out-Null
